$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two blank placeholder rows (rows 2 and 3), shifting the data
# (currently in rows 4-9) up to rows 2-7.
$ws.Rows("2:3").Delete()

# Restore the selection Excel leaves after a row delete.
$ws.Range("A2:XFD2").Select()
